$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as plain text, preserving default (unstyled) formatting,
# so numeric-looking strings (e.g. "0.125") are not auto-converted to numbers.
function Set-TextCell($cell, $value) {
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") '42.468.19'
$ws.Range("E2").Value = '  +0.71%  '

# Row 3
Set-TextCell $ws.Range("D3") '2.235.53'
$ws.Range("E3").Value = '  -0.52%  '

# Row 4
$ws.Range("E4").Value = '  +0.41%  '

# Row 5
Set-TextCell $ws.Range("D5") '244.18'
$ws.Range("E5").Value = '  -1.47%  '

# Row 6
Set-TextCell $ws.Range("D6") '0.628'
$ws.Range("E6").Value = '  +0.40%  '

# Row 7
Set-TextCell $ws.Range("D7") '74.65'
$ws.Range("E7").Value = '  -3.52%  '

# Row 8
$ws.Range("E8").Value = '  +0.12%  '

# Row 9
Set-TextCell $ws.Range("D9") '0.616'
$ws.Range("E9").Value = '  -2.82%  '

# Row 10
Set-TextCell $ws.Range("D10") '43.35'
$ws.Range("E10").Value = '  +3.29%  '

# Row 11
Set-TextCell $ws.Range("D11") '0.0948'
$ws.Range("E11").Value = '  -1.11%  '

# Row 12
Set-TextCell $ws.Range("D12") '7.12'
$ws.Range("E12").Value = '  -0.83%  '

# Row 13
$ws.Range("E13").Value = '  -0.11%  '

# Row 14
Set-TextCell $ws.Range("D14") '14.45'
$ws.Range("E14").Value = '  -2.79%  '

# Row 15
Set-TextCell $ws.Range("D15") '0.851'
$ws.Range("E15").Value = '  -1.30%  '

# Row 16
Set-TextCell $ws.Range("D16") '2.238.23'
$ws.Range("E16").Value = '  +0.80%  '

# Row 17
Set-TextCell $ws.Range("D17") '42.283.36'
$ws.Range("E17").Value = '  +0.62%  '

# Row 18
$ws.Range("E18").Value = '  +6.15%  '

# Row 19
Set-TextCell $ws.Range("D19") '6.15'
$ws.Range("E19").Value = '  +0.06%  '

# Row 20
Set-TextCell $ws.Range("D20") '71.91'
$ws.Range("E20").Value = '  -0.15%  '

# Row 21
Set-TextCell $ws.Range("D21") '10.16'
$ws.Range("E21").Value = '  +39.17%  '

# Row 22
Set-TextCell $ws.Range("D22") '230.59'
$ws.Range("E22").Value = '  -0.68%  '

# Row 23
$ws.Range("E23").Value = '  -6.44%  '

# Row 24
Set-TextCell $ws.Range("D24") '11.62'
$ws.Range("E24").Value = '  +1.97%  '

# Row 25
$ws.Range("E25").Value = '  +0.00%  '

# Row 26
Set-TextCell $ws.Range("D26") '3.64'
$ws.Range("E26").Value = '  -0.70%  '

# Row 27
$ws.Range("E27").Value = '  +0.08%  '

# Row 28
$ws.Range("E28").Value = '  +4.03%  '

# Row 29
Set-TextCell $ws.Range("D29") '166.61'
$ws.Range("E29").Value = '  -1.86%  '

# Row 30
Set-TextCell $ws.Range("D30") '20.81'
$ws.Range("E30").Value = '  +1.04%  '

# Row 31
$ws.Range("E31").Value = '  +20.36%  '

# Row 32
Set-TextCell $ws.Range("D32") '0.0809'
$ws.Range("E32").Value = '  -2.96%  '

# Row 33
Set-TextCell $ws.Range("D33") '0.118'
$ws.Range("E33").Value = '  -2.95%  '

# Row 34
Set-TextCell $ws.Range("B34") 'Stellar'
Set-TextCell $ws.Range("C34") 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell $ws.Range("D34") '0.125'
$ws.Range("E34").Value = '  -0.75%  '

# Row 35
Set-TextCell $ws.Range("B35") 'InjectiveProtocol'
Set-TextCell $ws.Range("C35") 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell $ws.Range("D35") '29.45'
$ws.Range("E35").Value = '  -13.44%  '

# Row 36
Set-TextCell $ws.Range("D36") '4.54'
$ws.Range("E36").Value = '  -0.18%  '

# Row 37
Set-TextCell $ws.Range("D37") '0.0310'
$ws.Range("E37").Value = '  +2.56%  '

# Row 38
Set-TextCell $ws.Range("D38") '13.19'
$ws.Range("E38").Value = '  -8.09%  '

# Row 39
Set-TextCell $ws.Range("D39") '2.16'
$ws.Range("E39").Value = '  -1.63%  '

# Row 40
Set-TextCell $ws.Range("D40") '5.68'
$ws.Range("E40").Value = '  -4.67%  '

# Row 41
Set-TextCell $ws.Range("D41") '63.23'
$ws.Range("E41").Value = '  +3.42%  '

# Row 42
Set-TextCell $ws.Range("D42") '0.200'
$ws.Range("E42").Value = '  -2.15%  '

# Row 43
Set-TextCell $ws.Range("B43") 'FraxShare'
Set-TextCell $ws.Range("C43") 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell $ws.Range("D43") '8.81'
$ws.Range("E43").Value = '  +1.02%  '

# Row 44
Set-TextCell $ws.Range("B44") 'Aave'
Set-TextCell $ws.Range("C44") 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws.Range("D44") '105.85'
$ws.Range("E44").Value = '  -6.64%  '

# Row 45
Set-TextCell $ws.Range("D45") '0.102'
$ws.Range("E45").Value = '  +1.68%  '

# Row 46
Set-TextCell $ws.Range("D46") '0.996'
$ws.Range("E46").Value = '  -0.19%  '

# Row 47
Set-TextCell $ws.Range("D47") '2.39'
$ws.Range("E47").Value = '  +4.55%  '

# Row 48
Set-TextCell $ws.Range("D48") '1.13'
$ws.Range("E48").Value = '  -0.76%  '

# Row 49
Set-TextCell $ws.Range("D49") '1.17'
$ws.Range("E49").Value = '  -0.16%  '

# Row 50
$ws.Range("E50").Value = '  +1.31%  '

# Row 51
Set-TextCell $ws.Range("D51") '4.10'
$ws.Range("E51").Value = '  -2.71%  '
